# "React - Learned JSX"
# Mark the "Completed" column (E) as "Yes" for the roadmap rows that were
# finished, on Sheet2 of the JS Roadmap workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Activate()

$completedRows = @(69, 81, 82, 89, 90, 98, 102, 103, 104, 114, 120, 123, 124, 135)

foreach ($r in $completedRows) {
    $ws.Range("E$r").Value = "Yes"
}

# Leave the view scrolled down to where editing stopped, same as Excel
# would persist in the saved sheetView / selection.
$ws.Range("E148").Select()
$excel.ActiveWindow.ScrollRow = 127
$excel.ActiveWindow.ScrollColumn = 1
